$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 0.6643128802000001
$ws.Range("R2").Value = 5.978815921800001
$ws.Range("S2").Value = 0.8271666313262851
$ws.Range("T2").Value = 0.8271666313262852

# Row 3
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("S3").Value = 0.09421438109281059
$ws.Range("T3").Value = 0.09421438109281059

# Row 4
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("S4").Value = 0.07861898758090437
$ws.Range("T4").Value = 0.07861898758090438
